# Actualiza base de datos del Estado de Cuenta:
# - Se eliminan los periodos de mora anteriores (2404/2506) y se agregan los nuevos (2507/2404)
# - Se actualizan los valores de mora y el salario basico asociado

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 16: nuevo periodo de mora 2507
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Fila 17: nuevo periodo de mora 2404
$ws.Range("E17").Value = "2404"
$ws.Range("F17").Value = 22533
$ws.Range("G17").Value = 1423500

# El contenido mas ancho recalcula el ajuste de las columnas (bestFit)
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
